$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-01 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-02 Thursday", 2)
$d.Content.Find.Execute("23×64=", $true, $false, $false, $false, $false, $true, 1, $false, "69×40=", 2)
$d.Content.Find.Execute("83×53=", $true, $false, $false, $false, $false, $true, 1, $false, "97×60=", 2)
$d.Content.Find.Execute("69×28=", $true, $false, $false, $false, $false, $true, 1, $false, "72×29=", 2)
$d.Content.Find.Execute("26×21=", $true, $false, $false, $false, $false, $true, 1, $false, "77×91=", 2)
$d.Content.Find.Execute("82×42=", $true, $false, $false, $false, $false, $true, 1, $false, "71×40=", 2)
$d.Content.Find.Execute("78×46=", $true, $false, $false, $false, $false, $true, 1, $false, "45×60=", 2)
$d.Content.Find.Execute("67×47=", $true, $false, $false, $false, $false, $true, 1, $false, "18×21=", 2)
$d.Content.Find.Execute("94×59=", $true, $false, $false, $false, $false, $true, 1, $false, "40×95=", 2)
$d.Content.Find.Execute("59×31=", $true, $false, $false, $false, $false, $true, 1, $false, "33×93=", 2)
$d.Content.Find.Execute("21×71=", $true, $false, $false, $false, $false, $true, 1, $false, "22×39=", 2)
$d.Content.Find.Execute("52×65=", $true, $false, $false, $false, $false, $true, 1, $false, "55×39=", 2)
$d.Content.Find.Execute("74×55=", $true, $false, $false, $false, $false, $true, 1, $false, "43×17=", 2)
$d.Content.Find.Execute("32×42=", $true, $false, $false, $false, $false, $true, 1, $false, "67×69=", 2)
$d.Content.Find.Execute("58×86=", $true, $false, $false, $false, $false, $true, 1, $false, "94×85=", 2)
$d.Content.Find.Execute("49×16=", $true, $false, $false, $false, $false, $true, 1, $false, "98×54=", 2)
$d.Content.Find.Execute("76×16=", $true, $false, $false, $false, $false, $true, 1, $false, "29×61=", 2)
$d.Content.Find.Execute("61×61=", $true, $false, $false, $false, $false, $true, 1, $false, "40×75=", 2)
$d.Content.Find.Execute("48×96=", $true, $false, $false, $false, $false, $true, 1, $false, "39×63=", 2)
$d.Content.Find.Execute("40×91=", $true, $false, $false, $false, $false, $true, 1, $false, "48×71=", 2)
$d.Content.Find.Execute("21×80=", $true, $false, $false, $false, $false, $true, 1, $false, "17×63=", 2)
$d.Content.Find.Execute("80×47=", $true, $false, $false, $false, $false, $true, 1, $false, "93×24=", 2)
$d.Content.Find.Execute("35×77=", $true, $false, $false, $false, $false, $true, 1, $false, "73×87=", 2)
$d.Content.Find.Execute("78×51=", $true, $false, $false, $false, $false, $true, 1, $false, "60×45=", 2)
$d.Content.Find.Execute("77×65=", $true, $false, $false, $false, $false, $true, 1, $false, "29×55=", 2)
$d.Content.Find.Execute("99×47=", $true, $false, $false, $false, $false, $true, 1, $false, "34×76=", 2)
